$d = $word.ActiveDocument

# Header date line
$d.Content.Find.Execute("2026-01-07 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-08 Thursday", 2)

# Table of division problems.
# Most values are globally unique in the document, so a simple Find/Replace
# across the whole content is safe. The two cells that originally read
# "78÷6=" are ambiguous for Find/Replace (they map to two different new
# values), so those two are addressed directly via the Tables collection
# using their (row, column) position instead.

$d.Content.Find.Execute("22÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷7=", 2)
$d.Content.Find.Execute("65÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "15÷7=", 2)
$d.Content.Find.Execute("64÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "27÷5=", 2)
$d.Content.Find.Execute("73÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "25÷9=", 2)
$d.Content.Find.Execute("12÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷6=", 2)

$d.Content.Find.Execute("68÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "35÷8=", 2)
$d.Content.Find.Execute("15÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "40÷7=", 2)
$d.Content.Find.Execute("95÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "98÷5=", 2)
$d.Content.Find.Execute("49÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷2=", 2)
$d.Content.Find.Execute("66÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷4=", 2)

$d.Content.Find.Execute("23÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "89÷6=", 2)
$d.Content.Find.Execute("56÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "79÷6=", 2)
$d.Content.Find.Execute("23÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "30÷4=", 2)
$d.Content.Find.Execute("24÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "49÷3=", 2)

# Row 3 (1-based row 9), last column: "78÷6=" -> "21÷6="
$d.Tables.Item(1).Cell(9, 5).Range.Text = "21÷6="

# Row 4 (1-based row 13), first column: "78÷6=" -> "11÷8="
$d.Tables.Item(1).Cell(13, 1).Range.Text = "11÷8="

$d.Content.Find.Execute("39÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "12÷2=", 2)
$d.Content.Find.Execute("40÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "37÷9=", 2)
$d.Content.Find.Execute("53÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "42÷7=", 2)
$d.Content.Find.Execute("46÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷5=", 2)

$d.Content.Find.Execute("60÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "47÷9=", 2)
$d.Content.Find.Execute("91÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷5=", 2)
$d.Content.Find.Execute("17÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "13÷7=", 2)
$d.Content.Find.Execute("92÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "58÷5=", 2)
$d.Content.Find.Execute("13÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "62÷4=", 2)
